$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells in row 1: strip the interior spaces from these
# column-title labels (e.g. "S1 Eps" -> "S1Eps", "Top 1 Word Count" -> "Top1WordCount").
$ws.Range("B1").Value = "S1Eps"
$ws.Range("C1").Value = "S1Lines"
$ws.Range("J1").Value = "S2Eps"
$ws.Range("K1").Value = "S2Lines"
$ws.Range("AJ1").Value = "S3Eps"
$ws.Range("AK1").Value = "S3Lines"
$ws.Range("BB1").Value = "S4Eps"
$ws.Range("BC1").Value = "S4Lines"
$ws.Range("BZ1").Value = "S5Eps"
$ws.Range("CA1").Value = "S5Lines"
$ws.Range("CX1").Value = "S6Eps"
$ws.Range("CY1").Value = "S6Lines"
$ws.Range("DT1").Value = "S7Eps"
$ws.Range("DU1").Value = "S7Lines"
$ws.Range("EH1").Value = "TotalEps"
$ws.Range("EI1").Value = "TotalLines"
$ws.Range("EJ1").Value = "Top1Line"
$ws.Range("EK1").Value = "Top1LineCount"
$ws.Range("EL1").Value = "Top2Line"
$ws.Range("EM1").Value = "Top2LineCount"
$ws.Range("EN1").Value = "Top3Line"
$ws.Range("EO1").Value = "Top3LineCount"
$ws.Range("EP1").Value = "Top1Word"
$ws.Range("EQ1").Value = "Top1WordCount"
$ws.Range("ER1").Value = "Top2Word"
$ws.Range("ES1").Value = "Top2WordCount"
$ws.Range("ET1").Value = "Top3Word"
$ws.Range("EU1").Value = "Top3WordCount"
$ws.Range("EV1").Value = "Top4Word"
$ws.Range("EW1").Value = "Top4WordCount"
$ws.Range("EX1").Value = "Top5Word"
$ws.Range("EY1").Value = "Top5WordCount"
$ws.Range("EZ1").Value = "Top6Word"
$ws.Range("FA1").Value = "Top6WordCount"
$ws.Range("FB1").Value = "Top7Word"
$ws.Range("FC1").Value = "Top7WordCount"
$ws.Range("FD1").Value = "Top8Word"
$ws.Range("FE1").Value = "Top8WordCount"
$ws.Range("FF1").Value = "Top9Word"
$ws.Range("FG1").Value = "Top9WordCount"
$ws.Range("FH1").Value = "Top10Word"
$ws.Range("FI1").Value = "Top10WordCount"

# View-state changes captured in the source workbook: new zoom level and
# new active selection (the user had scrolled/zoomed before saving).
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("FK16").Select()
